$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "CONNECTIVITY"
$wb.Worksheets.Item(2).Name = "COORDINATES"
$wb.Worksheets.Item(3).Name = "FREE NODES"
